$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$c = $ws.Range("B14")
Write-Host $c.HorizontalAlignment()
Write-Host $c.VerticalAlignment()
Write-Host $c.WrapText()
